$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 988.2406999999999
$ws.Cells.Item(17, 9).Value = 591.64
$ws.Cells.Item(17, 10).Value = 1330.138
$ws.Cells.Item(17, 11).Value = 1774.92
$ws.Cells.Item(17, 12).Value = 3990.414
$ws.Cells.Item(17, 13).Value = -1606.92
$ws.Cells.Item(17, 14).Value = -4326.414
$ws.Cells.Item(53, 8).Value = 182.38461
$ws.Cells.Item(53, 9).Value = 205.16667
$ws.Cells.Item(53, 10).Value = 162.85715
$ws.Cells.Item(53, 11).Value = 205.16667
$ws.Cells.Item(53, 12).Value = 162.85715
$ws.Cells.Item(53, 13).Value = 431.83333
$ws.Cells.Item(53, 14).Value = -1436.85715
$ws.Cells.Item(112, 8).Value = 5904.1934
$ws.Cells.Item(112, 10).Value = 6497.5
$ws.Cells.Item(112, 12).Value = 19492.5
$ws.Cells.Item(112, 14).Value = -21708.5
$ws.Cells.Item(116, 8).Value = 9080.200000000001
$ws.Cells.Item(116, 9).Value = 10683.667
$ws.Cells.Item(116, 11).Value = 10683.667
$ws.Cells.Item(116, 13).Value = -7241.666999999999
$ws.Cells.Item(132, 8).Value = 2091.459
$ws.Cells.Item(132, 9).Value = 2002.1818
$ws.Cells.Item(132, 10).Value = 2909.8333
$ws.Cells.Item(132, 11).Value = 6006.5454
$ws.Cells.Item(132, 12).Value = 8729.499899999999
$ws.Cells.Item(132, 13).Value = -3476.5454
$ws.Cells.Item(132, 14).Value = -13789.4999
$ws.Cells.Item(137, 8).Value = 1393.6061
$ws.Cells.Item(137, 9).Value = 1405.6316
$ws.Cells.Item(137, 10).Value = 1377.2858
$ws.Cells.Item(137, 11).Value = 4216.8948
$ws.Cells.Item(137, 12).Value = 4131.857400000001
$ws.Cells.Item(137, 13).Value = -1666.8948
$ws.Cells.Item(137, 14).Value = -9231.857400000001
$ws.Cells.Item(138, 8).Value = 1268.51
$ws.Cells.Item(138, 9).Value = 486.59573
$ws.Cells.Item(138, 10).Value = 1961.9056
$ws.Cells.Item(138, 11).Value = 1459.78719
$ws.Cells.Item(138, 12).Value = 5885.7168
$ws.Cells.Item(138, 13).Value = 3680.21281
$ws.Cells.Item(138, 14).Value = -16165.7168
$ws.Cells.Item(141, 8).Value = 2553.1167
$ws.Cells.Item(141, 9).Value = 784.81396
$ws.Cells.Item(141, 10).Value = 7025.8823
$ws.Cells.Item(141, 11).Value = 2354.44188
$ws.Cells.Item(141, 12).Value = 21077.6469
$ws.Cells.Item(141, 13).Value = 2825.55812
$ws.Cells.Item(141, 14).Value = -31437.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 840349.75
$ws.Cells.Item(32, 9).Value = 952254.6
$ws.Cells.Item(32, 10).Value = 19714
$ws.Cells.Item(32, 11).Value = 952254.6
$ws.Cells.Item(32, 12).Value = 19714
$ws.Cells.Item(32, 13).Value = -951967.6
$ws.Cells.Item(32, 14).Value = -20288
$ws.Cells.Item(61, 8).Value = 2493.9575
$ws.Cells.Item(61, 9).Value = 2083.6128
$ws.Cells.Item(61, 10).Value = 3289
$ws.Cells.Item(61, 11).Value = 2083.6128
$ws.Cells.Item(61, 12).Value = 3289
$ws.Cells.Item(61, 13).Value = -1871.6128
$ws.Cells.Item(61, 14).Value = -3713
$ws.Cells.Item(74, 8).Value = 856.3929000000001
$ws.Cells.Item(74, 9).Value = 616.5714
$ws.Cells.Item(74, 11).Value = 616.5714
$ws.Cells.Item(74, 13).Value = 257.4286
$ws.Cells.Item(77, 8).Value = 856.3929000000001
$ws.Cells.Item(77, 9).Value = 616.5714
$ws.Cells.Item(77, 11).Value = 3082.857
$ws.Cells.Item(77, 13).Value = 1285.143
$ws.Cells.Item(132, 8).Value = 3784.3655
$ws.Cells.Item(132, 9).Value = 3293.6316
$ws.Cells.Item(132, 10).Value = 5116.357
$ws.Cells.Item(132, 11).Value = 9880.8948
$ws.Cells.Item(132, 12).Value = 15349.071
$ws.Cells.Item(132, 13).Value = -7350.8948
$ws.Cells.Item(132, 14).Value = -20409.071
$ws.Cells.Item(136, 8).Value = 2493.9575
$ws.Cells.Item(136, 9).Value = 2083.6128
$ws.Cells.Item(136, 10).Value = 3289
$ws.Cells.Item(136, 11).Value = 6250.8384
$ws.Cells.Item(136, 12).Value = 9867
$ws.Cells.Item(136, 13).Value = -3700.8384
$ws.Cells.Item(136, 14).Value = -14967

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1970.3214
$ws.Cells.Item(20, 9).Value = 1626
$ws.Cells.Item(20, 10).Value = 2502.4546
$ws.Cells.Item(20, 11).Value = 1626
$ws.Cells.Item(20, 12).Value = 2502.4546
$ws.Cells.Item(20, 13).Value = -1379
$ws.Cells.Item(20, 14).Value = -2996.4546
$ws.Cells.Item(86, 8).Value = 4464.4287
$ws.Cells.Item(86, 9).Value = 4636.923
$ws.Cells.Item(86, 10).Value = 2222
$ws.Cells.Item(86, 11).Value = 4636.923
$ws.Cells.Item(86, 12).Value = 2222
$ws.Cells.Item(86, 13).Value = -3513.923
$ws.Cells.Item(86, 14).Value = -4468
$ws.Cells.Item(89, 8).Value = 4464.4287
$ws.Cells.Item(89, 9).Value = 4636.923
$ws.Cells.Item(89, 10).Value = 2222
$ws.Cells.Item(89, 11).Value = 23184.615
$ws.Cells.Item(89, 12).Value = 11110
$ws.Cells.Item(89, 13).Value = -17568.615
$ws.Cells.Item(89, 14).Value = -22342
$ws.Cells.Item(105, 8).Value = 13892166
$ws.Cells.Item(105, 9).Value = 13892166
$ws.Cells.Item(105, 11).Value = 13892166
$ws.Cells.Item(105, 13).Value = -13890419
$ws.Cells.Item(134, 8).Value = 2965.5789
$ws.Cells.Item(134, 9).Value = 2939.9092
$ws.Cells.Item(134, 10).Value = 3000.875
$ws.Cells.Item(134, 11).Value = 8819.7276
$ws.Cells.Item(134, 12).Value = 9002.625
$ws.Cells.Item(134, 13).Value = -6284.7276
$ws.Cells.Item(134, 14).Value = -14072.625
$ws.Cells.Item(135, 8).Value = 46645
$ws.Cells.Item(135, 10).Value = 46645
$ws.Cells.Item(135, 12).Value = 46645
$ws.Cells.Item(135, 14).Value = -56785

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4522.204
$ws.Cells.Item(31, 10).Value = 6504.8
$ws.Cells.Item(31, 12).Value = 6504.8
$ws.Cells.Item(31, 14).Value = -7094.8
$ws.Cells.Item(34, 8).Value = 4522.204
$ws.Cells.Item(34, 10).Value = 6504.8
$ws.Cells.Item(34, 12).Value = 6504.8
$ws.Cells.Item(34, 14).Value = -6908.8
$ws.Cells.Item(58, 8).Value = 1015.4375
$ws.Cells.Item(58, 9).Value = 755.7778
$ws.Cells.Item(58, 11).Value = 755.7778
$ws.Cells.Item(58, 13).Value = -552.7778
$ws.Cells.Item(122, 8).Value = 1937.3684
$ws.Cells.Item(122, 9).Value = 1461
$ws.Cells.Item(122, 10).Value = 1993.4117
$ws.Cells.Item(122, 11).Value = 4383
$ws.Cells.Item(122, 12).Value = 5980.2351
$ws.Cells.Item(122, 13).Value = -1933
$ws.Cells.Item(122, 14).Value = -10880.2351
$ws.Cells.Item(132, 8).Value = 3473557.8
$ws.Cells.Item(132, 9).Value = 1332.561
$ws.Cells.Item(132, 10).Value = 23810878
$ws.Cells.Item(132, 11).Value = 3997.683
$ws.Cells.Item(132, 12).Value = 71432634
$ws.Cells.Item(132, 13).Value = -1467.683
$ws.Cells.Item(132, 14).Value = -71437694
$ws.Cells.Item(134, 8).Value = 2687.3594
$ws.Cells.Item(134, 9).Value = 2844.426
$ws.Cells.Item(134, 10).Value = 1839.2
$ws.Cells.Item(134, 11).Value = 8533.278
$ws.Cells.Item(134, 12).Value = 5517.6
$ws.Cells.Item(134, 13).Value = -5998.278
$ws.Cells.Item(134, 14).Value = -10587.6
$ws.Cells.Item(136, 8).Value = 1015.4375
$ws.Cells.Item(136, 9).Value = 755.7778
$ws.Cells.Item(136, 11).Value = 2267.3334
$ws.Cells.Item(136, 13).Value = 282.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1141.1305
$ws.Cells.Item(5, 9).Value = 338.8421
$ws.Cells.Item(5, 10).Value = 1705.7037
$ws.Cells.Item(5, 11).Value = 1016.5263
$ws.Cells.Item(5, 12).Value = 5117.1111
$ws.Cells.Item(5, 13).Value = -904.5263
$ws.Cells.Item(5, 14).Value = -5341.1111
$ws.Cells.Item(131, 8).Value = 2698.9285
$ws.Cells.Item(131, 10).Value = 2979.9033
$ws.Cells.Item(131, 12).Value = 8939.7099
$ws.Cells.Item(131, 14).Value = -19019.7099
$ws.Cells.Item(135, 8).Value = 1141.1305
$ws.Cells.Item(135, 9).Value = 338.8421
$ws.Cells.Item(135, 10).Value = 1705.7037
$ws.Cells.Item(135, 11).Value = 3049.5789
$ws.Cells.Item(135, 12).Value = 15351.3333
$ws.Cells.Item(135, 13).Value = -514.5789
$ws.Cells.Item(135, 14).Value = -20421.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 677.5
$ws.Cells.Item(97, 9).Value = 739
$ws.Cells.Item(97, 10).Value = 542.2
$ws.Cells.Item(97, 11).Value = 739
$ws.Cells.Item(97, 12).Value = 542.2
$ws.Cells.Item(97, 13).Value = -243
$ws.Cells.Item(97, 14).Value = -1534.2
$ws.Cells.Item(107, 8).Value = 490.33334
$ws.Cells.Item(107, 9).Value = 472
$ws.Cells.Item(107, 10).Value = 499.5
$ws.Cells.Item(107, 11).Value = 472
$ws.Cells.Item(107, 12).Value = 499.5
$ws.Cells.Item(107, 13).Value = 1448
$ws.Cells.Item(107, 14).Value = -4339.5
$ws.Cells.Item(117, 8).Value = 55154.25
$ws.Cells.Item(117, 10).Value = 55154.25
$ws.Cells.Item(117, 12).Value = 55154.25
$ws.Cells.Item(117, 14).Value = -62038.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 2878.7273
$ws.Cells.Item(100, 9).Value = 2785.111
$ws.Cells.Item(100, 11).Value = 2785.111
$ws.Cells.Item(100, 13).Value = -2244.111
$ws.Cells.Item(132, 8).Value = 2172.0527
$ws.Cells.Item(132, 9).Value = 1949.05
$ws.Cells.Item(132, 10).Value = 2696.7646
$ws.Cells.Item(132, 11).Value = 5847.15
$ws.Cells.Item(132, 12).Value = 8090.293799999999
$ws.Cells.Item(132, 13).Value = -3317.15
$ws.Cells.Item(132, 14).Value = -13150.2938
$ws.Cells.Item(136, 8).Value = 3969668
$ws.Cells.Item(136, 9).Value = 1398.6875
$ws.Cells.Item(136, 10).Value = 16668130
$ws.Cells.Item(136, 11).Value = 4196.0625
$ws.Cells.Item(136, 12).Value = 50004390
$ws.Cells.Item(136, 13).Value = -1646.0625
$ws.Cells.Item(136, 14).Value = -50009490

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2724.3333
$ws.Cells.Item(122, 9).Value = 2679
$ws.Cells.Item(122, 10).Value = 3155
$ws.Cells.Item(122, 11).Value = 8037
$ws.Cells.Item(122, 12).Value = 9465
$ws.Cells.Item(122, 13).Value = -5587
$ws.Cells.Item(122, 14).Value = -14365
$ws.Cells.Item(132, 8).Value = 4505927.5
$ws.Cells.Item(132, 9).Value = 1620.75
$ws.Cells.Item(132, 10).Value = 12821570
$ws.Cells.Item(132, 11).Value = 4862.25
$ws.Cells.Item(132, 12).Value = 38464710
$ws.Cells.Item(132, 13).Value = -2332.25
$ws.Cells.Item(132, 14).Value = -38469770
$ws.Cells.Item(136, 8).Value = 1946.241
$ws.Cells.Item(136, 9).Value = 1834.295
$ws.Cells.Item(136, 10).Value = 2256.6365
$ws.Cells.Item(136, 11).Value = 5502.885
$ws.Cells.Item(136, 12).Value = 6769.9095
$ws.Cells.Item(136, 13).Value = -2952.885
$ws.Cells.Item(136, 14).Value = -11869.9095
